$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list data (price + 1h volume % changes)
$ws.Range("D2").Value = "30.220.09"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.858.06"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.58"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4704"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2888"
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06554"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.56"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07938"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.42"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.851.09"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.101"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6775"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.47"
$ws.Range("E16").Value = "  -3.94%  "
$ws.Range("D17").Value = "30.212.45"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.62"
$ws.Range("E18").Value = "  +7.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007634"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("D21").Value = "2.110.01"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.238"
$ws.Range("E23").Value = "  -4.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.150"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.12"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.149"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.85"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.937"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.395"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09849"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.468"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.299"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.996"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04692"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6985"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.613"
$ws.Range("E39").Value = "  +3.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.327"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.41"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.929"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8415"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4137"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.15"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.075"
$ws.Range("E47").Value = "  -1.45%  "

# Row 48 and 49 swap places (EnergySwap now ranks above Maker) with updated data
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.156"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "935.88"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.98"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05662"
$ws.Range("E51").Value = "  +0.68%  "
